$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 61731.793
$ws.Range("I132").Value = 68038.2
$ws.Range("J132").Value = 6550.75
$ws.Range("K132").Value = 204114.6
$ws.Range("L132").Value = 19652.25
$ws.Range("M132").Value = -201584.6
$ws.Range("N132").Value = -24712.25
$ws.Range("H137").Value = 904695.9
$ws.Range("I137").Value = 2673.4092
$ws.Range("K137").Value = 8020.2276
$ws.Range("M137").Value = -5470.2276
$ws.Range("H141").Value = 3000.0
$ws.Range("I141").Value = 3000.0
$ws.Range("J141").Value = 3000.0
$ws.Range("K141").Value = 9000.0
$ws.Range("L141").Value = 9000.0
$ws.Range("M141").Value = -3820.0
$ws.Range("N141").Value = -19360.0

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 359.25
$ws.Range("I5").Value = 254.8
$ws.Range("J5").Value = 533.3333
$ws.Range("K5").Value = 254.8
$ws.Range("L5").Value = 533.3333
$ws.Range("M5").Value = -142.8
$ws.Range("N5").Value = -757.3333
$ws.Range("H32").Value = 25417.467
$ws.Range("I32").Value = 25417.467
$ws.Range("K32").Value = 25417.467
$ws.Range("M32").Value = -25130.467
$ws.Range("H36").Value = 3742.0
$ws.Range("I36").Value = 4113.0
$ws.Range("J36").Value = 3000.0
$ws.Range("K36").Value = 4113.0
$ws.Range("L36").Value = 3000.0
$ws.Range("M36").Value = -3767.0
$ws.Range("N36").Value = -3692.0
$ws.Range("H51").Value = 45000.0
$ws.Range("J51").Value = 45000.0
$ws.Range("L51").Value = 45000.0
$ws.Range("N51").Value = -46512.0
$ws.Range("H61").Value = 3031386.2
$ws.Range("I61").Value = 3031386.2
$ws.Range("K61").Value = 3031386.2
$ws.Range("M61").Value = -3031174.2
$ws.Range("H74").Value = 2115.1155
$ws.Range("I74").Value = 908.6316
$ws.Range("J74").Value = 5389.857
$ws.Range("K74").Value = 908.6316
$ws.Range("L74").Value = 5389.857
$ws.Range("M74").Value = -34.63160000000005
$ws.Range("N74").Value = -7137.857
$ws.Range("H77").Value = 2115.1155
$ws.Range("I77").Value = 908.6316
$ws.Range("J77").Value = 5389.857
$ws.Range("K77").Value = 4543.158
$ws.Range("L77").Value = 26949.285
$ws.Range("M77").Value = -175.1580000000004
$ws.Range("N77").Value = -35685.285
$ws.Range("H136").Value = 3031386.2
$ws.Range("I136").Value = 3031386.2
$ws.Range("K136").Value = 9094158.600000001
$ws.Range("M136").Value = -9091608.600000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 359.25
$ws.Range("I4").Value = 254.8
$ws.Range("J4").Value = 533.3333
$ws.Range("K4").Value = 254.8
$ws.Range("L4").Value = 533.3333
$ws.Range("M4").Value = -139.8
$ws.Range("N4").Value = -763.3333
$ws.Range("H107").Value = 3489.5334
$ws.Range("J107").Value = 4812.6
$ws.Range("L107").Value = 4812.6
$ws.Range("N107").Value = -8652.6
$ws.Range("H134").Value = 581267.1
$ws.Range("I134").Value = 542912.1
$ws.Range("K134").Value = 1628736.3
$ws.Range("M134").Value = -1626201.3
$ws.Range("H140").Value = 100000.0
$ws.Range("J140").Value = 100000.0
$ws.Range("L140").Value = 100000.0
$ws.Range("N140").Value = -110360.0

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10400.706
$ws.Range("I31").Value = 3750.7646
$ws.Range("K31").Value = 3750.7646
$ws.Range("M31").Value = -3455.7646
$ws.Range("H34").Value = 10400.706
$ws.Range("I34").Value = 3750.7646
$ws.Range("K34").Value = 3750.7646
$ws.Range("M34").Value = -3548.7646
$ws.Range("H58").Value = 691079.0
$ws.Range("I58").Value = 951493.3
$ws.Range("K58").Value = 951493.3
$ws.Range("M58").Value = -951290.3
$ws.Range("H59").Value = 97212.8
$ws.Range("I59").Value = 1500.0
$ws.Range("K59").Value = 1500.0
$ws.Range("M59").Value = -355.0
$ws.Range("H74").Value = 80313.5
$ws.Range("J74").Value = 80313.5
$ws.Range("L74").Value = 80313.5
$ws.Range("N74").Value = -82061.5
$ws.Range("H77").Value = 80313.5
$ws.Range("J77").Value = 80313.5
$ws.Range("L77").Value = 240940.5
$ws.Range("N77").Value = -249676.5
$ws.Range("H97").Value = 27999.5
$ws.Range("J97").Value = 27999.5
$ws.Range("L97").Value = 27999.5
$ws.Range("N97").Value = -29981.5
$ws.Range("H132").Value = 10227.762
$ws.Range("I132").Value = 2466.6155
$ws.Range("K132").Value = 7399.8465
$ws.Range("M132").Value = -4869.8465
$ws.Range("H136").Value = 691079.0
$ws.Range("I136").Value = 951493.3
$ws.Range("K136").Value = 2854479.9
$ws.Range("M136").Value = -2851929.9

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 0.0
$ws.Range("J69").Value = 0.0
$ws.Range("L69").Value = 0.0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0.0
$ws.Range("J72").Value = 0.0
$ws.Range("L72").Value = 0.0
$ws.Range("N72").ClearContents()
$ws.Range("H121").Value = 28572588.0
$ws.Range("I121").Value = 200000000.0
$ws.Range("J121").Value = 1351.5
$ws.Range("K121").Value = 600000000.0
$ws.Range("L121").Value = 4054.5
$ws.Range("M121").Value = -599998690.0
$ws.Range("N121").Value = -6674.5
$ws.Range("H129").Value = 3995.2856
$ws.Range("I129").Value = 0.0
$ws.Range("J129").Value = 3995.2856
$ws.Range("K129").Value = 0.0
$ws.Range("L129").Value = 11985.8568
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -21985.8568

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 83900.0
$ws.Range("J39").Value = 83900.0
$ws.Range("L39").Value = 83900.0
$ws.Range("N39").Value = -84964.0
$ws.Range("H80").Value = 192598.81
$ws.Range("I80").Value = 258493.0
$ws.Range("J80").Value = 4329.7144
$ws.Range("K80").Value = 258493.0
$ws.Range("L80").Value = 4329.7144
$ws.Range("M80").Value = -257495.0
$ws.Range("N80").Value = -6325.7144
$ws.Range("H83").Value = 192598.81
$ws.Range("I83").Value = 258493.0
$ws.Range("J83").Value = 4329.7144
$ws.Range("K83").Value = 1292465.0
$ws.Range("L83").Value = 21648.572
$ws.Range("M83").Value = -1287473.0
$ws.Range("N83").Value = -31632.572
$ws.Range("H113").Value = 9498.25
$ws.Range("I113").Value = 0.0
$ws.Range("K113").Value = 0.0
$ws.Range("M113").ClearContents()
$ws.Range("H126").Value = 880050.8
$ws.Range("I126").Value = 1517808.1
$ws.Range("K126").Value = 4553424.300000001
$ws.Range("M126").Value = -4550954.300000001
$ws.Range("H132").Value = 507354.03
$ws.Range("I132").Value = 578222.94
$ws.Range("K132").Value = 1734668.82
$ws.Range("M132").Value = -1732138.82

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3911.7222
$ws.Range("I7").Value = 3672.3572
$ws.Range("K7").Value = 3672.3572
$ws.Range("M7").Value = -3560.3572
$ws.Range("H68").Value = 2654.5386
$ws.Range("I68").Value = 2082.1667
$ws.Range("K68").Value = 2082.1667
$ws.Range("M68").Value = -1333.1667
$ws.Range("H71").Value = 2654.5386
$ws.Range("I71").Value = 2082.1667
$ws.Range("K71").Value = 10410.8335
$ws.Range("M71").Value = -6666.833500000001
$ws.Range("H76").Value = 8525.333
$ws.Range("J76").Value = 8525.333
$ws.Range("L76").Value = 8525.333
$ws.Range("N76").Value = -9201.333
$ws.Range("H79").Value = 8525.333
$ws.Range("J79").Value = 8525.333
$ws.Range("L79").Value = 8525.333
$ws.Range("N79").Value = -10865.333
$ws.Range("H98").Value = 90083.4
$ws.Range("J98").Value = 90083.4
$ws.Range("L98").Value = 90083.4
$ws.Range("N98").Value = -96073.4
$ws.Range("H100").Value = 10415.77
$ws.Range("I100").Value = 2733.4443
$ws.Range("J100").Value = 27701.0
$ws.Range("K100").Value = 2733.4443
$ws.Range("L100").Value = 27701.0
$ws.Range("M100").Value = -2192.4443
$ws.Range("N100").Value = -28783.0
$ws.Range("H126").Value = 3911.7222
$ws.Range("I126").Value = 3672.3572
$ws.Range("K126").Value = 11017.0716
$ws.Range("M126").Value = -8547.0716
$ws.Range("H136").Value = 3728.1904
$ws.Range("I136").Value = 2886.4666
$ws.Range("K136").Value = 8659.399800000001
$ws.Range("M136").Value = -6109.399800000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 29633.334
$ws.Range("I3").Value = 29000.0
$ws.Range("J3").Value = 29950.0
$ws.Range("K3").Value = 29000.0
$ws.Range("L3").Value = 29950.0
$ws.Range("M3").Value = -28886.0
$ws.Range("N3").Value = -30178.0
$ws.Range("H51").Value = 14999.5
$ws.Range("I51").Value = 14999.5
$ws.Range("K51").Value = 14999.5
$ws.Range("M51").Value = -14489.5
$ws.Range("H107").Value = 2043.85
$ws.Range("I107").Value = 1083.5385
$ws.Range("J107").Value = 3827.2856
$ws.Range("K107").Value = 3250.6155
$ws.Range("L107").Value = 11481.8568
$ws.Range("M107").Value = -1330.6155
$ws.Range("N107").Value = -15321.8568
$ws.Range("H122").Value = 2228.4736
$ws.Range("I122").Value = 1875.0625
$ws.Range("J122").Value = 4113.3335
$ws.Range("K122").Value = 5625.1875
$ws.Range("L122").Value = 12340.0005
$ws.Range("M122").Value = -3175.1875
$ws.Range("N122").Value = -17240.0005
